# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) is re-sorted from descending
# (1903 -> 1809) to ascending (1809 -> 1903) order, and the "Valor Mora"
# amounts in column F follow the same re-sort (the values that belonged to
# periods 1903 and 1809 swap rows accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the period labels into ascending order
$ws.Range("E16").Value = "1809"
$ws.Range("E17").Value = "1810"
$ws.Range("E18").Value = "1811"
$ws.Range("E19").Value = "1812"
$ws.Range("E20").Value = "1901"
$ws.Range("E21").Value = "1902"
$ws.Range("E22").Value = "1903"

# Carry the "Valor Mora" values along with the re-sorted periods
# (only the rows for period 1903 and 1809 actually change value)
$ws.Range("F16").Value = 31249
$ws.Range("F22").Value = 26041
